$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.401.06"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "1.960.68"
$ws.Range("E3").Value = "  -4.14%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "249.34"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "0.602"
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "53.51"
$ws.Range("E8").Value = "  -10.51%  "
$ws.Range("D9").Value = "0.368"
$ws.Range("E9").Value = "  -6.76%  "
$ws.Range("D10").Value = "0.0747"
$ws.Range("E10").Value = "  -7.79%  "
$ws.Range("D11").Value = "0.100"
$ws.Range("E11").Value = "  -3.84%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.265.28"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "13.90"
$ws.Range("E13").Value = "  -9.50%  "
$ws.Range("D14").Value = "20.87"
$ws.Range("E14").Value = "  -5.52%  "
$ws.Range("D15").Value = "0.763"
$ws.Range("E15").Value = "  -10.89%  "
$ws.Range("D16").Value = "5.07"
$ws.Range("E16").Value = "  -7.42%  "
$ws.Range("D17").Value = "1.975.52"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "36.315.53"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D19").Value = "68.72"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("D20").Value = "0.0₃0805"
$ws.Range("E20").Value = "  -6.78%  "
$ws.Range("D21").Value = "229.72"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "4.95"
$ws.Range("E22").Value = "  -6.33%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "2.47"
$ws.Range("E24").Value = "  -3.65%  "
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "162.42"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").Value = "8.66"
$ws.Range("E27").Value = "  -8.19%  "
$ws.Range("D28").Value = "18.98"
$ws.Range("E28").Value = "  -5.03%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.121"
$ws.Range("E29").Value = "  -12.15%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.30"
$ws.Range("E30").Value = "  -5.66%  "
$ws.Range("D31").Value = "0.117"
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("D32").Value = "4.41"
$ws.Range("E32").Value = "  -8.26%  "
$ws.Range("D33").Value = "0.0618"
$ws.Range("E33").Value = "  -9.68%  "
$ws.Range("D34").Value = "4.25"
$ws.Range("E34").Value = "  -5.75%  "
$ws.Range("D35").Value = "2.30"
$ws.Range("E35").Value = "  -10.06%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "1.81"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("D39").Value = "5.37"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").Value = "2.99"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.15"
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.417.59"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("D43").Value = "0.0894"
$ws.Range("E43").Value = "  -8.84%  "
$ws.Range("D44").Value = "0.0203"
$ws.Range("E44").Value = "  -6.70%  "
$ws.Range("D45").Value = "86.76"
$ws.Range("E45").Value = "  -6.00%  "
$ws.Range("D46").Value = "15.13"
$ws.Range("E46").Value = "  -9.79%  "
$ws.Range("D47").Value = "0.992"
$ws.Range("E47").Value = "  -7.27%  "
$ws.Range("D48").Value = "2.85"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("D49").Value = "6.75"
$ws.Range("E49").Value = "  -10.52%  "
$ws.Range("D50").Value = "2.160.51"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("D51").Value = "1.88"
$ws.Range("E51").Value = "  -11.77%  "
